$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update country ordering: Angola moves ahead of Lituania (new entrant),
#    Lituania and Cabo Verde shift down one row, Mauritania keeps its row.
$ws.Cells.Item(117,1).Value = "Angola"
$ws.Cells.Item(118,1).Value = "Lituania"
$ws.Cells.Item(119,1).Value = "Cabo Verde"
$ws.Cells.Item(120,1).Value = "Mauritania"

# 2. Update statistic values for the affected rows
$ws.Cells.Item(4,2).Value = 8514514
$ws.Cells.Item(4,3).Value = 57861
$ws.Cells.Item(4,4).Value = 5540114
$ws.Cells.Item(4,5).Value = 2748353
$ws.Cells.Item(4,7).Value = 825
$ws.Cells.Item(4,8).Value = 226047

$ws.Cells.Item(6,4).Value = 4721593
$ws.Cells.Item(6,5).Value = 397524

$ws.Cells.Item(10,2).Value = 974139
$ws.Cells.Item(10,3).Value = 8256
$ws.Cells.Item(10,4).Value = 876731
$ws.Cells.Item(10,5).Value = 68136
$ws.Cells.Item(10,7).Value = 170
$ws.Cells.Item(10,8).Value = 29272

$ws.Cells.Item(12,2).Value = 874118
$ws.Cells.Item(12,3).Value = 3242
$ws.Cells.Item(12,4).Value = 788494
$ws.Cells.Item(12,5).Value = 51749
$ws.Cells.Item(12,7).Value = 55
$ws.Cells.Item(12,8).Value = 33875

$ws.Cells.Item(39,2).Value = 129944
$ws.Cells.Item(39,3).Value = 273
$ws.Cells.Item(39,4).Value = 126866
$ws.Cells.Item(39,5).Value = 2854

$ws.Cells.Item(47,2).Value = 105705
$ws.Cells.Item(47,3).Value = 158
$ws.Cells.Item(47,4).Value = 98413
$ws.Cells.Item(47,5).Value = 1150
$ws.Cells.Item(47,7).Value = 12
$ws.Cells.Item(47,8).Value = 6142

$ws.Cells.Item(92,2).Value = 21570
$ws.Cells.Item(92,3).Value = 64
$ws.Cells.Item(92,5).Value = 1028
$ws.Cells.Item(92,7).Value = 1
$ws.Cells.Item(92,8).Value = 425

$ws.Cells.Item(97,2).Value = 16772
$ws.Cells.Item(97,3).Value = 169
$ws.Cells.Item(97,5).Value = 4631

$ws.Cells.Item(109,2).Value = 11066
$ws.Cells.Item(109,3).Value = 14
$ws.Cells.Item(109,4).Value = 10362
$ws.Cells.Item(109,5).Value = 401

$ws.Cells.Item(117,2).Value = 8049
$ws.Cells.Item(117,3).Value = 220
$ws.Cells.Item(117,4).Value = 3037
$ws.Cells.Item(117,5).Value = 4761
$ws.Cells.Item(117,7).Value = 3
$ws.Cells.Item(117,8).Value = 251

$ws.Cells.Item(118,2).Value = 7928
$ws.Cells.Item(118,3).Value = 202
$ws.Cells.Item(118,4).Value = 3276
$ws.Cells.Item(118,5).Value = 4534
$ws.Cells.Item(118,7).Value = 5
$ws.Cells.Item(118,8).Value = 118

$ws.Cells.Item(119,2).Value = 7901
$ws.Cells.Item(119,3).Value = 101
$ws.Cells.Item(119,4).Value = 6792
$ws.Cells.Item(119,5).Value = 1022
$ws.Cells.Item(119,8).Value = 87

$ws.Cells.Item(120,2).Value = 7634
$ws.Cells.Item(120,3).Value = 13
$ws.Cells.Item(120,4).Value = 7359
$ws.Cells.Item(120,5).Value = 112

# 3. Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Octubre de 2020 a las 01:03"
